$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-04-20 Saturday" "2024-04-21 Sunday"
Replace-Text "23÷8=2, 7" "58÷9=6, 4"
Replace-Text "71÷2=35, 1" "28÷8=3, 4"
Replace-Text "25÷5=5, 0" "59÷4=14, 3"
Replace-Text "33÷8=4, 1" "47÷5=9, 2"
Replace-Text "54÷7=7, 5" "80÷5=16, 0"
Replace-Text "60÷4=15, 0" "75÷9=8, 3"
Replace-Text "63÷3=21, 0" "41÷4=10, 1"
Replace-Text "41÷9=4, 5" "39÷4=9, 3"
Replace-Text "89÷5=17, 4" "89÷6=14, 5"
Replace-Text "43÷3=14, 1" "54÷8=6, 6"
Replace-Text "50÷3=16, 2" "79÷6=13, 1"
Replace-Text "86÷4=21, 2" "44÷5=8, 4"
Replace-Text "55÷3=18, 1" "42÷6=7, 0"
Replace-Text "18÷8=2, 2" "37÷2=18, 1"
Replace-Text "22÷4=5, 2" "78÷5=15, 3"
Replace-Text "46÷3=15, 1" "85÷8=10, 5"
Replace-Text "93÷6=15, 3" "63÷8=7, 7"
Replace-Text "14÷6=2, 2" "49÷3=16, 1"
Replace-Text "24÷5=4, 4" "71÷9=7, 8"
Replace-Text "56÷6=9, 2" "89÷7=12, 5"
Replace-Text "38÷5=7, 3" "86÷7=12, 2"
Replace-Text "26÷9=2, 8" "31÷6=5, 1"
Replace-Text "99÷8=12, 3" "55÷3=18, 1"
Replace-Text "73÷7=10, 3" "31÷2=15, 1"
Replace-Text "44÷9=4, 8" "12÷8=1, 4"
